$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark all TC (test case) Runmode values to "Y"
$ws.Range("C2").Value = "Y"
$ws.Range("C3").Value = "Y"
